$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Golang Architect / Principal Backend Architect || Atlanta, GA (Onsite)"
$ws.Range("B16").Value = "https://www.dice.com/job-detail/4d57826e-1249-42b5-a805-223a1887a5db"
$ws.Range("C16").Value = "Atlanta, Georgia"
$ws.Range("D16").Value = "Contract, Third Party"
$ws.Range("E16").Value = "Depends on Experience"
$ws.Range("F16").Value = "Galactic Minds Inc."
